$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New NPC rows (29-32), continuing the existing table of NPCs.
# Columns: A=id, B=name, C=real_name, D=type, E=game_map_id, I=x_position, J=y_position

$rows = @(
    @{ A = 28; B = "WanderingCleric"; C = "Wandering Cleric"; D = 2; E = "Delusional Memories"; I = 112; J = 32 },
    @{ A = 29; B = "JesterofTime";    C = "Jester of Time";    D = 2; E = "Delusional Memories"; I = 448; J = 320 },
    @{ A = 30; B = "HolyKnight";      C = "Holy Knight";       D = 2; E = "Delusional Memories"; I = 16;  J = 224 },
    @{ A = 31; B = "FaithlessChild";  C = "Faithless Child";   D = 2; E = "Delusional Memories"; I = 96;  J = 432 }
)

$startRow = 29
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J
}

# The new "Delusional Memories" map name is longer than anything previously in
# column E (game_map_id), so its best-fit width grows to fit the new text.
$ws.Columns.Item(5).ColumnWidth = 22.6
